$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 5 hours (5/24 of a day) of time spent in "week 3" (column E) for Ben (row 6) and Vincent (row 7)
$ws.Range("E6").Value = 0.30208333333333331
$ws.Range("E7").Value = 0.27083333333333331

# Update the active cell selection to F8 (as recorded in the workbook after the edit)
$ws.Range("F8").Select()
